$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 1.327089666666667
$ws.Cells.Item(2,8).Value = 3.981269
$ws.Cells.Item(2,9).Value = 0.0007725509783306236
$ws.Cells.Item(2,10).Value = 0.0007725509783306238
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 2.773245666666666
$ws.Cells.Item(2,14).Value = 8.319737
$ws.Cells.Item(2,15).Value = 0.1610278523700304
$ws.Cells.Item(2,16).Value = 0.1610278523700304
$ws.Cells.Item(2,17).Value = 3.680345667361444
$ws.Cells.Item(2,18).Value = 33.123111006253
$ws.Cells.Item(2,19).Value = 0.0001244022248869462
$ws.Cells.Item(2,20).Value = 0.0001244022248869462

# Row 3
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 1.327089666666667
$ws.Cells.Item(3,8).Value = 3.981269
$ws.Cells.Item(3,9).Value = 0.0007725509783306236
$ws.Cells.Item(3,10).Value = 0.0007725509783306238
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 11.68452533333333
$ws.Cells.Item(3,14).Value = 35.053576
$ws.Cells.Item(3,15).Value = 0.6784591942232839
$ws.Cells.Item(3,16).Value = 0.6784591942232838
$ws.Cells.Item(3,17).Value = 15.50641282977155
$ws.Cells.Item(3,18).Value = 139.557715467944
$ws.Cells.Item(3,19).Value = 0.0005241443142546046
$ws.Cells.Item(3,20).Value = 0.0005241443142546046

# Row 4
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 1.327089666666667
$ws.Cells.Item(4,8).Value = 3.981269
$ws.Cells.Item(4,9).Value = 0.0007725509783306236
$ws.Cells.Item(4,10).Value = 0.0007725509783306238
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 2.764378
$ws.Cells.Item(4,14).Value = 8.293134
$ws.Cells.Item(4,15).Value = 0.1605129534066858
$ws.Cells.Item(4,16).Value = 0.1605129534066858
$ws.Cells.Item(4,17).Value = 3.668577478560667
$ws.Cells.Item(4,18).Value = 33.017197307046
$ws.Cells.Item(4,19).Value = 0.0001240044391890729
$ws.Cells.Item(4,20).Value = 0.0001240044391890729

# Row 5
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 1678.019409
$ws.Cells.Item(5,8).Value = 5034.058227
$ws.Cells.Item(5,9).Value = 0.9768409540380654
$ws.Cells.Item(5,10).Value = 0.9768409540380655
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 2.773245666666666
$ws.Cells.Item(5,14).Value = 8.319737
$ws.Cells.Item(5,15).Value = 0.1610278523700304
$ws.Cells.Item(5,16).Value = 0.1610278523700304
$ws.Cells.Item(5,17).Value = 4653.560054591811
$ws.Cells.Item(5,18).Value = 41882.0404913263
$ws.Cells.Item(5,19).Value = 0.1572986009358412
$ws.Cells.Item(5,20).Value = 0.1572986009358412

# Row 6
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 1678.019409
$ws.Cells.Item(6,8).Value = 5034.058227
$ws.Cells.Item(6,9).Value = 0.9768409540380654
$ws.Cells.Item(6,10).Value = 0.9768409540380655
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 11.68452533333333
$ws.Cells.Item(6,14).Value = 35.053576
$ws.Cells.Item(6,15).Value = 0.6784591942232839
$ws.Cells.Item(6,16).Value = 0.6784591942232838
$ws.Cells.Item(6,17).Value = 19606.86029428553
$ws.Cells.Item(6,18).Value = 176461.7426485697
$ws.Cells.Item(6,19).Value = 0.6627467265609698
$ws.Cells.Item(6,20).Value = 0.6627467265609698

# Row 7
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 1678.019409
$ws.Cells.Item(7,8).Value = 5034.058227
$ws.Cells.Item(7,9).Value = 0.9768409540380654
$ws.Cells.Item(7,10).Value = 0.9768409540380655
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 2.764378
$ws.Cells.Item(7,14).Value = 8.293134
$ws.Cells.Item(7,15).Value = 0.1605129534066858
$ws.Cells.Item(7,16).Value = 0.1605129534066858
$ws.Cells.Item(7,17).Value = 4638.679937812602
$ws.Cells.Item(7,18).Value = 41748.11944031341
$ws.Cells.Item(7,19).Value = 0.1567956265412545
$ws.Cells.Item(7,20).Value = 0.1567956265412545

# Row 8
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 38.45556733333333
$ws.Cells.Item(8,8).Value = 115.366702
$ws.Cells.Item(8,9).Value = 0.02238649498360385
$ws.Cells.Item(8,10).Value = 0.02238649498360385
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 2.773245666666666
$ws.Cells.Item(8,14).Value = 8.319737
$ws.Cells.Item(8,15).Value = 0.1610278523700304
$ws.Cells.Item(8,16).Value = 0.1610278523700304
$ws.Cells.Item(8,17).Value = 106.6467354663749
$ws.Cells.Item(8,18).Value = 959.8206191973741
$ws.Cells.Item(8,19).Value = 0.003604849209302187
$ws.Cells.Item(8,20).Value = 0.003604849209302187

# Row 9
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 38.45556733333333
$ws.Cells.Item(9,8).Value = 115.366702
$ws.Cells.Item(9,9).Value = 0.02238649498360385
$ws.Cells.Item(9,10).Value = 0.02238649498360385
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 11.68452533333333
$ws.Cells.Item(9,14).Value = 35.053576
$ws.Cells.Item(9,15).Value = 0.6784591942232839
$ws.Cells.Item(9,16).Value = 0.6784591942232838
$ws.Cells.Item(9,17).Value = 449.3350507140391
$ws.Cells.Item(9,18).Value = 4044.015456426352
$ws.Cells.Item(9,19).Value = 0.01518832334805946
$ws.Cells.Item(9,20).Value = 0.01518832334805946

# Row 10
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 38.45556733333333
$ws.Cells.Item(10,8).Value = 115.366702
$ws.Cells.Item(10,9).Value = 0.02238649498360385
$ws.Cells.Item(10,10).Value = 0.02238649498360385
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 2.764378
$ws.Cells.Item(10,14).Value = 8.293134
$ws.Cells.Item(10,15).Value = 0.1605129534066858
$ws.Cells.Item(10,16).Value = 0.1605129534066858
$ws.Cells.Item(10,17).Value = 106.3057243137853
$ws.Cells.Item(10,18).Value = 956.7515188240681
$ws.Cells.Item(10,19).Value = 0.003593322426242211
$ws.Cells.Item(10,20).Value = 0.00359332242624221
